# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Price cells that look numeric are prefixed with a leading apostrophe so
# Excel keeps them as literal text (matching the sheet's existing text
# cells) instead of silently coercing e.g. "111.90" into the number 111.9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.557.01"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "2.637.93"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'111.90"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'325.92"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").Value = "'39.53"
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("D11").Value = "'20.09"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "'0.0810"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "'7.50"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Value = "3.051.40"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "2.640.05"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "'0.854"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").Value = "49.513.57"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'13.11"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "'268.64"
$ws.Range("E23").Value = "  -2.92%  "
$ws.Range("D24").Value = "'69.11"
$ws.Range("E24").Value = "  -4.09%  "
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").Value = "'26.09"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'10.20"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").Value = "'34.68"
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("D32").Value = "'49.56"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").Value = "'5.49"
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "'19.03"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("E37").Value = "  +4.58%  "
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").Value = "'128.44"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("D41").Value = "'22.94"
$ws.Range("E41").Value = "  +3.18%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'0.0330"
$ws.Range("E44").Value = "  +5.11%  "
$ws.Range("D45").Value = "2.057.23"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").Value = "'3.27"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").Value = "'2.14"
$ws.Range("E47").Value = "  +8.33%  "
$ws.Range("E48").Value = "  -5.88%  "
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("D51").Value = "'58.66"
$ws.Range("E51").Value = "  -0.06%  "
